$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.736.62"
$ws.Range("E2").Value = "  +0.95%  "

# Row 3
$ws.Range("D3").Value = "1.889.95"
$ws.Range("E3").Value = "  +1.03%  "

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").Value = "'247.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.28%  "

# Row 6
$ws.Range("E6").Value = "  +0.13%  "

# Row 7
$ws.Range("D7").Value = "'0.4737"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("D8").Value = "'0.2932"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.40%  "

# Row 9
$ws.Range("D9").Value = "'0.06534"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.55%  "

# Row 10
$ws.Range("D10").Value = "'22.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.20%  "

# Row 11
$ws.Range("D11").Value = "'0.07796"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.99%  "

# Row 12
$ws.Range("D12").Value = "'97.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.51%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.892.75"
$ws.Range("E13").Value = "  +1.16%  "

# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.7409"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.16%  "

# Row 15
$ws.Range("D15").Value = "'5.254"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.07%  "

# Row 16
$ws.Range("D16").Value = "'286.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.65%  "

# Row 17
$ws.Range("D17").Value = "30.719.42"
$ws.Range("E17").Value = "  +0.99%  "

# Row 18
$ws.Range("D18").Value = "'13.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.03%  "

# Row 19
$ws.Range("D19").Value = "'0.000007536"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.15%  "

# Row 20
$ws.Range("D20").Value = "'1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.13%  "

# Row 21
$ws.Range("D21").Value = "2.141.76"
$ws.Range("E21").Value = "  +1.09%  "

# Row 22
$ws.Range("D22").Value = "'5.321"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.62%  "

# Row 23
$ws.Range("D23").Value = "'1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.25%  "

# Row 24
$ws.Range("D24").Value = "'6.272"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.49%  "

# Row 25
$ws.Range("D25").Value = "'9.223"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.81%  "

# Row 26
$ws.Range("D26").Value = "'164.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.76%  "

# Row 27
$ws.Range("D27").Value = "'19.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.06%  "

# Row 28
$ws.Range("D28").Value = "'1.920"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.25%  "

# Row 29
$ws.Range("D29").Value = "'1.346"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.48%  "

# Row 30
$ws.Range("D30").Value = "'0.09775"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.55%  "

# Row 31
$ws.Range("D31").Value = "'1.490"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.99%  "

# Row 32
$ws.Range("D32").Value = "'4.306"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.61%  "

# Row 33
$ws.Range("D33").Value = "'4.190"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.47%  "

# Row 34
$ws.Range("D34").Value = "'0.04895"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.25%  "

# Row 35
$ws.Range("D35").Value = "'1.129"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.76%  "

# Row 36
$ws.Range("D36").Value = "'0.7002"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.66%  "

# Row 37
$ws.Range("D37").Value = "'2.730"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.44%  "

# Row 38
$ws.Range("D38").Value = "'0.01900"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.62%  "

# Row 39
$ws.Range("D39").Value = "'2.839"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.64%  "

# Row 40
$ws.Range("D40").Value = "'76.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.68%  "

# Row 41
$ws.Range("D41").Value = "'6.327"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.21%  "

# Row 42
$ws.Range("D42").Value = "'2.010"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.01%  "

# Row 43
$ws.Range("D43").Value = "'0.4287"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.20%  "

# Row 44
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "'1.002"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.17%  "

# Row 45
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'0.8393"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.11%  "

# Row 46
$ws.Range("D46").Value = "'101.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.03%  "

# Row 47
$ws.Range("D47").Value = "'9.596"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.96%  "

# Row 48
$ws.Range("D48").Value = "'7.038"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.57%  "

# Row 49
$ws.Range("D49").Value = "'35.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.15%  "

# Row 50
$ws.Range("D50").Value = "'911.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.23%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05769"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.30%  "
